$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# ---------------------------------------------------------------------
# 1) Fill in missing data on existing rows 36 & 37 ("thomsonreuters"
#    answer column, and the missing G37 email + hyperlink).
# ---------------------------------------------------------------------
$ws.Range("L36").Value2 = "thomsonreuters"
$ws.Range("L37").Value2 = "thomsonreuters"

$ws.Range("G37").Value2 = "PLAUtestuser13@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("G37"), "mailto:PLAUtestuser13@mailinator.com") | Out-Null

# B37 no longer carries a hyperlink (text itself is unchanged)
$ws.Range("B37").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2) Append four new user rows (38-41), cloning the formatting of the
#    existing row 37 and then overwriting the content per row.
# ---------------------------------------------------------------------
$ws.Range("A37:L37").Copy($ws.Range("A38:L38"))
$ws.Range("A37:L37").Copy($ws.Range("A39:L39"))
$ws.Range("A37:L37").Copy($ws.Range("A40:L40"))
$ws.Range("A37:L37").Copy($ws.Range("A41:L41"))

# Row 38 - PLAUtestuser1
$ws.Range("A38").Value2 = "PLAUtestuser1"
$ws.Range("B38").Value2 = "P@ssword2"
$ws.Range("C38").Value2 = ""
$ws.Range("D38").Value2 = ""
$ws.Range("E38").Value2 = ""
$ws.Range("F38").Value2 = ""
$ws.Range("G38").Value2 = "PLAUtestuser1@mailinator.com"
$ws.Range("H38").Value2 = ""
$ws.Range("I38").Value2 = ""
$ws.Range("J38").Value2 = "ProdAUtestuser13"
$ws.Range("K38").Value2 = "ProdAUtestuser13@mailinator.com"
$ws.Range("L38").Value2 = ""

# Row 39 - PLAUtestuser2
$ws.Range("A39").Value2 = "PLAUtestuser2"
$ws.Range("B39").Value2 = "P@ssword2"
$ws.Range("C39").Value2 = ""
$ws.Range("D39").Value2 = ""
$ws.Range("E39").Value2 = ""
$ws.Range("F39").Value2 = ""
$ws.Range("G39").Value2 = "PLAUtestuser2@mailinator.com"
$ws.Range("H39").Value2 = ""
$ws.Range("I39").Value2 = ""
$ws.Range("J39").Value2 = "ProdAUtestuser13"
$ws.Range("K39").Value2 = "ProdAUtestuser13@mailinator.com"
$ws.Range("L39").Value2 = ""

# Row 40 - PLAUtestuser3
$ws.Range("A40").Value2 = "PLAUtestuser3"
$ws.Range("B40").Value2 = "P@ssword4"
$ws.Range("C40").Value2 = ""
$ws.Range("D40").Value2 = ""
$ws.Range("E40").Value2 = ""
$ws.Range("F40").Value2 = ""
$ws.Range("G40").Value2 = "PLAUtestuser3@mailinator.com"
$ws.Range("H40").Value2 = ""
$ws.Range("I40").Value2 = ""
$ws.Range("J40").Value2 = "ProdAUtestuser3"
$ws.Range("K40").Value2 = "ProdAUtestuser3@mailinator.com"
$ws.Range("L40").Value2 = "thomsonreuters"

# Row 41 - PLAUtestuser4
$ws.Range("A41").Value2 = "PLAUtestuser4"
$ws.Range("B41").Value2 = "P@ssword4"
$ws.Range("C41").Value2 = ""
$ws.Range("D41").Value2 = ""
$ws.Range("E41").Value2 = ""
$ws.Range("F41").Value2 = ""
$ws.Range("G41").Value2 = "PLAUtestuser4@mailinator.com"
$ws.Range("H41").Value2 = ""
$ws.Range("I41").Value2 = ""
$ws.Range("J41").Value2 = "ProdAUtestuser3"
$ws.Range("K41").Value2 = "ProdAUtestuser3@mailinator.com"
$ws.Range("L41").Value2 = "thomsonreuters"

# ---------------------------------------------------------------------
# 3) Hyperlinks for the new rows.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B38"), "mailto:P@ssword2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B39"), "mailto:P@ssword2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B40"), "mailto:P@ssword4") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B41"), "mailto:P@ssword4") | Out-Null

$ws.Hyperlinks.Add($ws.Range("G38"), "mailto:PLAUtestuser1@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G39"), "mailto:PLAUtestuser2@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G40"), "mailto:PLAUtestuser3@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G41"), "mailto:PLAUtestuser4@mailinator.com") | Out-Null

$ws.Hyperlinks.Add($ws.Range("K38"), "mailto:ProdAUtestuser13@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K39"), "mailto:ProdAUtestuser13@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K40"), "mailto:ProdAUtestuser3@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K41"), "mailto:ProdAUtestuser3@mailinator.com") | Out-Null

# ---------------------------------------------------------------------
# 4) Sheet view bookkeeping (scroll position / active selection).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C43").Select()
$excel.ActiveWindow.ScrollRow = 23

Write-Host "Edit complete"
